$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 25641968
$ws.Range("I40").Value = 41667412
$ws.Range("J40").Value = 1260
$ws.Range("K40").Value = 41667412
$ws.Range("L40").Value = 1260
$ws.Range("M40").Value = -41667237
$ws.Range("N40").Value = -1610
$ws.Range("H70").Value = 1008.1
$ws.Range("I70").Value = 620.25
$ws.Range("J70").Value = 1266.6666
$ws.Range("K70").Value = 1860.75
$ws.Range("L70").Value = 3799.9998
$ws.Range("M70").Value = -1590.75
$ws.Range("N70").Value = -4339.9998
$ws.Range("H73").Value = 1008.1
$ws.Range("I73").Value = 620.25
$ws.Range("J73").Value = 1266.6666
$ws.Range("K73").Value = 1860.75
$ws.Range("L73").Value = 3799.9998
$ws.Range("M73").Value = -924.75
$ws.Range("N73").Value = -5671.9998
$ws.Range("H92").Value = 825.6875
$ws.Range("I92").Value = 819.3333
$ws.Range("K92").Value = 819.3333
$ws.Range("M92").Value = 428.6667
$ws.Range("H98").Value = 3372.162
$ws.Range("I98").Value = 2994
$ws.Range("J98").Value = 9990
$ws.Range("K98").Value = 2994
$ws.Range("L98").Value = 9990
$ws.Range("M98").Value = -1496
$ws.Range("N98").Value = -12986
$ws.Range("H111").Value = 3774.5
$ws.Range("I111").Value = 4527.6665
$ws.Range("J111").Value = 3021.3333
$ws.Range("K111").Value = 13582.9995
$ws.Range("L111").Value = 9063.999899999999
$ws.Range("M111").Value = -10515.9995
$ws.Range("N111").Value = -15197.9999
$ws.Range("H122").Value = 3372.162
$ws.Range("I122").Value = 2994
$ws.Range("J122").Value = 9990
$ws.Range("K122").Value = 8982
$ws.Range("L122").Value = 29970
$ws.Range("M122").Value = -6532
$ws.Range("N122").Value = -34870
$ws.Range("H129").Value = 909.1539
$ws.Range("I129").Value = 750
$ws.Range("J129").Value = 938.0909
$ws.Range("K129").Value = 2250
$ws.Range("L129").Value = 2814.2727
$ws.Range("M129").Value = 2750
$ws.Range("N129").Value = -12814.2727
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2863.311
$ws.Range("I61").Value = 2483.111
$ws.Range("J61").Value = 3433.611
$ws.Range("K61").Value = 2483.111
$ws.Range("L61").Value = 3433.611
$ws.Range("M61").Value = -2271.111
$ws.Range("N61").Value = -3857.611
$ws.Range("H74").Value = 6561.143
$ws.Range("I74").Value = 1309.6875
$ws.Range("J74").Value = 13563.083
$ws.Range("K74").Value = 1309.6875
$ws.Range("L74").Value = 13563.083
$ws.Range("M74").Value = -435.6875
$ws.Range("N74").Value = -15311.083
$ws.Range("H77").Value = 6561.143
$ws.Range("I77").Value = 1309.6875
$ws.Range("J77").Value = 13563.083
$ws.Range("K77").Value = 6548.4375
$ws.Range("L77").Value = 67815.41500000001
$ws.Range("M77").Value = -2180.4375
$ws.Range("N77").Value = -76551.41500000001
$ws.Range("H132").Value = 84960320
$ws.Range("I132").Value = 123486430
$ws.Range("K132").Value = 370459290
$ws.Range("M132").Value = -370456760
$ws.Range("H136").Value = 2863.311
$ws.Range("I136").Value = 2483.111
$ws.Range("J136").Value = 3433.611
$ws.Range("K136").Value = 7449.333
$ws.Range("L136").Value = 10300.833
$ws.Range("M136").Value = -4899.333
$ws.Range("N136").Value = -15400.833
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1287.2903
$ws.Range("I86").Value = 1252.1666
$ws.Range("J86").Value = 1407.7142
$ws.Range("K86").Value = 1252.1666
$ws.Range("L86").Value = 1407.7142
$ws.Range("M86").Value = -129.1666
$ws.Range("N86").Value = -3653.7142
$ws.Range("H89").Value = 1287.2903
$ws.Range("I89").Value = 1252.1666
$ws.Range("J89").Value = 1407.7142
$ws.Range("K89").Value = 6260.833000000001
$ws.Range("L89").Value = 7038.571
$ws.Range("M89").Value = -644.8330000000005
$ws.Range("N89").Value = -18270.571
$ws.Range("H102").Value = 13518.333
$ws.Range("I102").Value = 13518.333
$ws.Range("K102").Value = 13518.333
$ws.Range("M102").Value = -10273.333
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3625968
$ws.Range("I31").Value = 1404.5714
$ws.Range("J31").Value = 5211714.5
$ws.Range("K31").Value = 1404.5714
$ws.Range("L31").Value = 5211714.5
$ws.Range("M31").Value = -1109.5714
$ws.Range("N31").Value = -5212304.5
$ws.Range("H34").Value = 3625968
$ws.Range("I34").Value = 1404.5714
$ws.Range("J34").Value = 5211714.5
$ws.Range("K34").Value = 1404.5714
$ws.Range("L34").Value = 5211714.5
$ws.Range("M34").Value = -1202.5714
$ws.Range("N34").Value = -5212118.5
$ws.Range("H50").Value = 9959.5
$ws.Range("J50").Value = 9959.5
$ws.Range("L50").Value = 9959.5
$ws.Range("N50").Value = -11209.5
$ws.Range("H59").Value = 15822.75
$ws.Range("J59").Value = 15797.429
$ws.Range("L59").Value = 15797.429
$ws.Range("N59").Value = -18087.429
$ws.Range("H68").Value = 17824
$ws.Range("J68").Value = 18219.111
$ws.Range("L68").Value = 18219.111
$ws.Range("N68").Value = -19717.111
$ws.Range("H71").Value = 17824
$ws.Range("J71").Value = 18219.111
$ws.Range("L71").Value = 54657.333
$ws.Range("N71").Value = -62145.333
$ws.Range("H74").Value = 14352.272
$ws.Range("I74").Value = 4692.5
$ws.Range("K74").Value = 4692.5
$ws.Range("M74").Value = -3818.5
$ws.Range("H77").Value = 14352.272
$ws.Range("I77").Value = 4692.5
$ws.Range("K77").Value = 14077.5
$ws.Range("M77").Value = -9709.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 1723.3334
$ws.Range("J22").Value = 1723.3334
$ws.Range("L22").Value = 5170.0002
$ws.Range("N22").Value = -5508.0002
$ws.Range("H27").Value = 1723.3334
$ws.Range("J27").Value = 1723.3334
$ws.Range("L27").Value = 5170.0002
$ws.Range("N27").Value = -5374.0002
$ws.Range("H56").Value = 3526.4
$ws.Range("I56").Value = 3526.4
$ws.Range("K56").Value = 3526.4
$ws.Range("M56").Value = -2996.4
$ws.Range("H119").Value = 812.375
$ws.Range("I119").Value = 812.375
$ws.Range("J119").Value = 0
$ws.Range("K119").Value = 2437.125
$ws.Range("L119").Value = 0
$ws.Range("M119").Value = 2400.875
$ws.Range("N119").ClearContents()
$ws.Range("H131").Value = 2290.2917
$ws.Range("I131").Value = 15307.25
$ws.Range("J131").Value = 1524.5883
$ws.Range("K131").Value = 45921.75
$ws.Range("L131").Value = 4573.7649
$ws.Range("M131").Value = -40881.75
$ws.Range("N131").Value = -14653.7649
$ws.Range("H132").Value = 47620390
$ws.Range("I132").Value = 66667590
$ws.Range("J132").Value = 2399.6667
$ws.Range("K132").Value = 600008310
$ws.Range("L132").Value = 21597.0003
$ws.Range("M132").Value = -600005780
$ws.Range("N132").Value = -26657.0003
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 59088320
$ws.Range("I132").Value = 102565976
$ws.Range("J132").Value = 2567367
$ws.Range("K132").Value = 307697928
$ws.Range("L132").Value = 7702101
$ws.Range("M132").Value = -307695398
$ws.Range("N132").Value = -7707161
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1825.7222
$ws.Range("I7").Value = 1708.3077
$ws.Range("J7").Value = 2131
$ws.Range("K7").Value = 1708.3077
$ws.Range("L7").Value = 2131
$ws.Range("M7").Value = -1596.3077
$ws.Range("N7").Value = -2355
$ws.Range("H16").Value = 2569.647
$ws.Range("I16").Value = 2569.647
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 2569.647
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -2399.647
$ws.Range("N16").ClearContents()
$ws.Range("H40").Value = 1749.5
$ws.Range("I40").Value = 1649.4
$ws.Range("J40").Value = 2250
$ws.Range("K40").Value = 1649.4
$ws.Range("L40").Value = 2250
$ws.Range("M40").Value = -1513.4
$ws.Range("N40").Value = -2522
$ws.Range("H46").Value = 1068.421
$ws.Range("I46").Value = 992.8570999999999
$ws.Range("J46").Value = 1280
$ws.Range("K46").Value = 992.8570999999999
$ws.Range("L46").Value = 1280
$ws.Range("M46").Value = -804.8570999999999
$ws.Range("N46").Value = -1656
$ws.Range("H122").Value = 2074.9167
$ws.Range("I122").Value = 1946.7715
$ws.Range("J122").Value = 2419.923
$ws.Range("K122").Value = 5840.3145
$ws.Range("L122").Value = 7259.768999999999
$ws.Range("M122").Value = -3390.3145
$ws.Range("N122").Value = -12159.769
$ws.Range("H126").Value = 1825.7222
$ws.Range("I126").Value = 1708.3077
$ws.Range("J126").Value = 2131
$ws.Range("K126").Value = 5124.9231
$ws.Range("L126").Value = 6393
$ws.Range("M126").Value = -2654.9231
$ws.Range("N126").Value = -11333
$ws.Range("H132").Value = 5088.8335
$ws.Range("I132").Value = 5433.1665
$ws.Range("J132").Value = 3367.1667
$ws.Range("K132").Value = 16299.4995
$ws.Range("L132").Value = 10101.5001
$ws.Range("M132").Value = -13769.4995
$ws.Range("N132").Value = -15161.5001
